$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-8 (name, King Checkout, King Stayover, Queen Checkout, Queen Stayover)
$data = @(
    @("Alexander", 0, 7, 4, 6),
    @("Andrea A", 1, 10, 2, 4),
    @("Johana", 2, 11, 3, 2),
    @("Julio", 1, 4, 2, 11),
    @("Mariana", 2, 9, 1, 4),
    @("Nestor", 2, 8, 1, 8),
    @("Tameka", 3, 3, 1, 9)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}

# Remove the old row 9 (Shae), shrinking the used range to A1:E8
$ws.Rows.Item(9).Delete()
